$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules Implementation")
$lo = $ws.ListObjects.Item("TabImp")
$col = $lo.ListColumns.Item("Situation Number")
$col.Name = "Situation ID"
